# Bugfixed the naive forecaster component module
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear stale forecast values in rows 2-6 (columns C and E)
$ws.Range("C2:C6").ClearContents()
$ws.Range("E2:E6").ClearContents()

# Updated forecast values for rows 7-19
$values = @{
    7  = @{ C = 1.097054137926201;  E = 1.37755776875883 }
    8  = @{ C = 1.385527545913412;  E = 1.329814931661888 }
    9  = @{ C = 1.296301936385214;  E = 1.355477993452414 }
    10 = @{ C = 2.441628883342295;  E = 1.386547975635688 }
    11 = @{ C = 2.565764046666463;  E = 1.833587970352424 }
    12 = @{ C = 1.263447557103259;  E = 1.485511920344451 }
    13 = @{ C = 2.117022522597423;  E = 1.745834498329324 }
    14 = @{ C = 2.149400276001101;  E = 1.76475225558832 }
    15 = @{ C = 2.453568910971748;  E = 2.131436976903012 }
    16 = @{ C = 0.812682184439506;  E = 1.556352278772266 }
    17 = @{ C = 0.9940067218177528; E = 1.820779918499094 }
    18 = @{ C = 1.634555928116921;  E = 1.554016159863814 }
    19 = @{ C = 0.6231570351797;    E = 1.581524829939718 }
}

foreach ($row in $values.Keys) {
    $ws.Range("C$row").Value = $values[$row].C
    $ws.Range("E$row").Value = $values[$row].E
}
